{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// 1) \"Ativa\u00e7\u00e3o: 01/01/2020\" -> \"Ativa\u00e7\u00e3o: 01/01/2025\"\n{\n  const oldText = \"Ativa\u00e7\u00e3o: 01/01/2020\";\n  const newText = \"Ativa\u00e7\u00e3o: 01/01/2025\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Append sentence to the Portuguese \"Programa\" paragraph\n{\n  const oldText = \"INTRODU\u00c7\u00c3O. Conceitos B\u00e1sicos. O perfil de solo. Defini\u00e7\u00e3o e nota\u00e7\u00e3o de horizontes e camadas. FORMA\u00c7\u00c3O DO SOLO. Fatores e processos de forma\u00e7\u00e3o. Intemperismo. ATRIBUTOS F\u00cdSICOS DO SOLO. Composi\u00e7\u00e3o volum\u00e9trica, granulometria e textura, estrutura e agrega\u00e7\u00e3o, cor, porosidade, densidade e compacta\u00e7\u00e3o, consist\u00eancia. CLASSIFICA\u00c7\u00c3O DO SOLO. Sistema brasileiro de classifica\u00e7\u00e3o de solos. Principais atributos morfol\u00f3gicos. Principais Classes de Solos. \u00c1GUA DO SOLO. Conceito e import\u00e2ncia. Constantes de umidade. Potencial total da \u00e1gua do solo e seus componentes. Curva caracter\u00edstica da \u00e1gua do solo. Movimento da \u00c1gua e de solutos no Solo. Aula pr\u00e1tica de campo: Descri\u00e7\u00e3o de perfil no campo. Aula pr\u00e1tica de laborat\u00f3rio: Caracteriza\u00e7\u00e3o e m\u00e9todos de determina\u00e7\u00e3o de atributos f\u00edsicos e h\u00eddricos do solo.\";\n  const newText = \"INTRODU\u00c7\u00c3O. Conceitos B\u00e1sicos. O perfil de solo. Defini\u00e7\u00e3o e nota\u00e7\u00e3o de horizontes e camadas. FORMA\u00c7\u00c3O DO SOLO. Fatores e processos de forma\u00e7\u00e3o. Intemperismo. ATRIBUTOS F\u00cdSICOS DO SOLO. Composi\u00e7\u00e3o volum\u00e9trica, granulometria e textura, estrutura e agrega\u00e7\u00e3o, cor, porosidade, densidade e compacta\u00e7\u00e3o, consist\u00eancia. CLASSIFICA\u00c7\u00c3O DO SOLO. Sistema brasileiro de classifica\u00e7\u00e3o de solos. Principais atributos morfol\u00f3gicos. Principais Classes de Solos. \u00c1GUA DO SOLO. Conceito e import\u00e2ncia. Constantes de umidade. Potencial total da \u00e1gua do solo e seus componentes. Curva caracter\u00edstica da \u00e1gua do solo. Movimento da \u00c1gua e de solutos no Solo. Aula pr\u00e1tica de campo: Descri\u00e7\u00e3o de perfil no campo. Aula pr\u00e1tica de laborat\u00f3rio: Caracteriza\u00e7\u00e3o e m\u00e9todos de determina\u00e7\u00e3o de atributos f\u00edsicos e h\u00eddricos do solo. A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Append sentence to the English (italic) \"Programa\" paragraph\n{\n  const oldText = \"INTRODUCTION. Basic Concepts. The soil profile. Horizons and layers definition and notation. SOIL FORMATION. Formation and factors processes. Weathering. PHYSICAL ATTRIBUTES OF THE SOIL. Volumetric composition, grain size and texture, structure and aggregation, color, porosity, density and compression, consistency. SOIL CLASSIFICATION. Brazilian system of soil classification. Main morphological attributes. Top Soil Classes. SOIL WATER. Concept and importance. Moisture constants. Total potential of Soil water and its components. Characteristic curve of soil water. Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil.\";\n  const newText = \"INTRODUCTION. Basic Concepts. The soil profile. Horizons and layers definition and notation. SOIL FORMATION. Formation and factors processes. Weathering. PHYSICAL ATTRIBUTES OF THE SOIL. Volumetric composition, grain size and texture, structure and aggregation, color, porosity, density and compression, consistency. SOIL CLASSIFICATION. Brazilian system of soil classification. Main morphological attributes. Top Soil Classes. SOIL WATER. Concept and importance. Moisture constants. Total potential of Soil water and its components. Characteristic curve of soil water. Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil. The discipline may have didactic trips to complement the content of the discipline.\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) Replace the \"M\u00e9todo:\" evaluation text\n{\n  const oldText = \"A avalia\u00e7\u00e3o ser\u00e1 feita mediante duas avalia\u00e7\u00f5es escritas de igual peso (P1 e P2). Alternativamente, essas avalia\u00e7\u00f5es escritas poder\u00e3o ser substitu\u00eddas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 ser\u00e3o dadas pela m\u00e9dia entre atividades desenvolvidas em aula, trabalhos e relat\u00f3rios de aulas pr\u00e1ticas.\";\n  const newText = \"O aluno poder\u00e1 optar por um dos dois crit\u00e9rios de avalia\u00e7\u00e3o para a NF (nota final).  Crit\u00e9rio 1: NF = m\u00e9dia obtida em todas atividades desenvolvidas, trabalhos e relat\u00f3rios ao longo do semestre. Crit\u00e9rio 2 (alternativo): NF = (P1+P2)/2, sendo P1 e P2 avalia\u00e7\u00f5es escritas individuais.\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 5) Replace the \"Crit\u00e9rio:\" evaluation text\n{\n  const oldText = \"O aluno poder\u00e1 optar por dois crit\u00e9rios de avalia\u00e7\u00e3o:Crit\u00e9rio 1: NF = (P1+P2)/2; ouCrit\u00e9rio 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avalia\u00e7\u00f5es escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relat\u00f3rios de aulas pr\u00e1ticas.\";\n  const newText = \"Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequ\u00eancia superior a 70%. Ser\u00e1 considerado aprovado o aluno que tenha obtido M\u00e9dia Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 6) Replace the \"Norma de recupera\u00e7\u00e3o:\" text (6,5 -> 5,0)\n{\n  const oldText = \"Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequ\u00eancia superior a 70%. Ser\u00e1 considerado aprovado o aluno que tenha obtido M\u00e9dia Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.\";\n  const newText = \"Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequ\u00eancia superior a 70%. Ser\u00e1 considerado aprovado o aluno que tenha obtido M\u00e9dia Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update activation date\n$findText = 'Ativa\u00e7\u00e3o: 01/01/2020'\n$replaceText = 'Ativa\u00e7\u00e3o: 01/01/2025'\n$range = $d.Content\n$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\nif (-not $found) {\n    throw 'Could not find text for: 1) Update activation date'\n}\n\n# 2) Append sentence to Portuguese Programa paragraph\n$findText = 'INTRODU\u00c7\u00c3O. Conceitos B\u00e1sicos. O perfil de solo. Defini\u00e7\u00e3o e nota\u00e7\u00e3o de horizontes e camadas. FORMA\u00c7\u00c3O DO SOLO. Fatores e processos de forma\u00e7\u00e3o. Intemperismo. ATRIBUTOS F\u00cdSICOS DO SOLO. Composi\u00e7\u00e3o volum\u00e9trica, granulometria e textura, estrutura e agrega\u00e7\u00e3o, cor, porosidade, densidade e compacta\u00e7\u00e3o, consist\u00eancia. CLASSIFICA\u00c7\u00c3O DO SOLO. Sistema brasileiro de classifica\u00e7\u00e3o de solos. Principais atributos morfol\u00f3gicos. Principais Classes de Solos. \u00c1GUA DO SOLO. Conceito e import\u00e2ncia. Constantes de umidade. Potencial total da \u00e1gua do solo e seus componentes. Curva caracter\u00edstica da \u00e1gua do solo. Movimento da \u00c1gua e de solutos no Solo. Aula pr\u00e1tica de campo: Descri\u00e7\u00e3o de perfil no campo. Aula pr\u00e1tica de laborat\u00f3rio: Caracteriza\u00e7\u00e3o e m\u00e9todos de determina\u00e7\u00e3o de atributos f\u00edsicos e h\u00eddricos do solo.'\n$replaceText = 'INTRODU\u00c7\u00c3O. Conceitos B\u00e1sicos. O perfil de solo. Defini\u00e7\u00e3o e nota\u00e7\u00e3o de horizontes e camadas. FORMA\u00c7\u00c3O DO SOLO. Fatores e processos de forma\u00e7\u00e3o. Intemperismo. ATRIBUTOS F\u00cdSICOS DO SOLO. Composi\u00e7\u00e3o volum\u00e9trica, granulometria e textura, estrutura e agrega\u00e7\u00e3o, cor, porosidade, densidade e compacta\u00e7\u00e3o, consist\u00eancia. CLASSIFICA\u00c7\u00c3O DO SOLO. Sistema brasileiro de classifica\u00e7\u00e3o de solos. Principais atributos morfol\u00f3gicos. Principais Classes de Solos. \u00c1GUA DO SOLO. Conceito e import\u00e2ncia. Constantes de umidade. Potencial total da \u00e1gua do solo e seus componentes. Curva caracter\u00edstica da \u00e1gua do solo. Movimento da \u00c1gua e de solutos no Solo. Aula pr\u00e1tica de campo: Descri\u00e7\u00e3o de perfil no campo. Aula pr\u00e1tica de laborat\u00f3rio: Caracteriza\u00e7\u00e3o e m\u00e9todos de determina\u00e7\u00e3o de atributos f\u00edsicos e h\u00eddricos do solo. A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.'\n$range = $d.Content\n$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\nif (-not $found) {\n    throw 'Could not find text for: 2) Append sentence to Portuguese Programa paragraph'\n}\n\n# 3) Append sentence to English (italic) Programa paragraph\n$findText = 'INTRODUCTION. Basic Concepts. The soil profile. Horizons and layers definition and notation. SOIL FORMATION. Formation and factors processes. Weathering. PHYSICAL ATTRIBUTES OF THE SOIL. Volumetric composition, grain size and texture, structure and aggregation, color, porosity, density and compression, consistency. SOIL CLASSIFICATION. Brazilian system of soil classification. Main morphological attributes. Top Soil Classes. SOIL WATER. Concept and importance. Moisture constants. Total potential of Soil water and its components. Characteristic curve of soil water. Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil.'\n$replaceText = 'INTRODUCTION. Basic Concepts. The soil profile. Horizons and layers definition and notation. SOIL FORMATION. Formation and factors processes. Weathering. PHYSICAL ATTRIBUTES OF THE SOIL. Volumetric composition, grain size and texture, structure and aggregation, color, porosity, density and compression, consistency. SOIL CLASSIFICATION. Brazilian system of soil classification. Main morphological attributes. Top Soil Classes. SOIL WATER. Concept and importance. Moisture constants. Total potential of Soil water and its components. Characteristic curve of soil water. Movement of water and solutes in soil. Field class practice: Profile description in the field. Laboratory class practice: Characterization and determination methods of physical and hydraulic properties of the soil. The discipline may have didactic trips to complement the content of the discipline.'\n$range = $d.Content\n$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\nif (-not $found) {\n    throw 'Could not find text for: 3) Append sentence to English (italic) Programa paragraph'\n}\n\n# 4) Replace Metodo evaluation text\n$findText = 'A avalia\u00e7\u00e3o ser\u00e1 feita mediante duas avalia\u00e7\u00f5es escritas de igual peso (P1 e P2). Alternativamente, essas avalia\u00e7\u00f5es escritas poder\u00e3o ser substitu\u00eddas por duas notas de igual peso (NOTA 1 e NOTA 2). Essas NOTAS 1 e 2 ser\u00e3o dadas pela m\u00e9dia entre atividades desenvolvidas em aula, trabalhos e relat\u00f3rios de aulas pr\u00e1ticas.'\n$replaceText = 'O aluno poder\u00e1 optar por um dos dois crit\u00e9rios de avalia\u00e7\u00e3o para a NF (nota final).  Crit\u00e9rio 1: NF = m\u00e9dia obtida em todas atividades desenvolvidas, trabalhos e relat\u00f3rios ao longo do semestre. Crit\u00e9rio 2 (alternativo): NF = (P1+P2)/2, sendo P1 e P2 avalia\u00e7\u00f5es escritas individuais.'\n$range = $d.Content\n$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\nif (-not $found) {\n    throw 'Could not find text for: 4) Replace Metodo evaluation text'\n}\n\n# 5) Replace Criterio evaluation text\n$findText = 'O aluno poder\u00e1 optar por dois crit\u00e9rios de avalia\u00e7\u00e3o:Crit\u00e9rio 1: NF = (P1+P2)/2; ouCrit\u00e9rio 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avalia\u00e7\u00f5es escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relat\u00f3rios de aulas pr\u00e1ticas.'\n$replaceText = 'Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequ\u00eancia superior a 70%. Ser\u00e1 considerado aprovado o aluno que tenha obtido M\u00e9dia Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.'\n$range = $d.Content\n$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\nif (-not $found) {\n    throw 'Could not find text for: 5) Replace Criterio evaluation text'\n}\n\n# 6) Replace Norma de recuperacao text (6,5 -> 5,0)\n$findText = 'Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 6,5 e frequ\u00eancia superior a 70%. Ser\u00e1 considerado aprovado o aluno que tenha obtido M\u00e9dia Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.'\n$replaceText = 'Exame Final (EF) para alunos com Nota Final (NF) maior ou igual a 3,0 e menor do que 5,0 e frequ\u00eancia superior a 70%. Ser\u00e1 considerado aprovado o aluno que tenha obtido M\u00e9dia Final (MF) igual ou maior do que 5,0, sendo MF = (NF+EF)/2.'\n$range = $d.Content\n$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\nif (-not $found) {\n    throw 'Could not find text for: 6) Replace Norma de recuperacao text (6,5 -> 5,0)'\n}\n"}
